{"js": "// Replace the empty curly-quote placeholders (\"\u201c\u201d\") throughout the\n// document with the real name of the botanical garden, and tidy up the\n// last paragraph (which previously held the stray \"_GoBack\" bookmark\n// used only to mark where the placeholder text used to be typed).\n\nconst body = context.document.body;\n\n// --- 1) Title paragraph -----------------------------------------------\n// \"...jard\u00edn bot\u00e1nico \"\u201d.\u201d  ->  \"...jard\u00edn bot\u00e1nico de la Universidad del\n// Mar campus Puerto Escondido.\"\nconst titleMatches = body.search(\"\u201c\u201d\", { matchCase: false });\ntitleMatches.load(\"text\");\nawait context.sync();\nif (titleMatches.items.length > 0) {\n  titleMatches.items[0].insertText(\n    \"de la Universidad del Mar campus Puerto Escondido\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- 2) \"Objetivo general\" paragraph ------------------------------------\n// \"...visitantes del jard\u00edn bot\u00e1nico \"\u201d de la Universidad del Mar.\u201d  ->\n// \"...visitantes del jard\u00edn bot\u00e1nico de la Universidad del Mar campus\n// Puerto Escondido.\"\nconst generalMatches = body.search(\"\u201c\u201d de la Universidad del Mar\", {\n  matchCase: false,\n});\ngeneralMatches.load(\"text\");\nawait context.sync();\nif (generalMatches.items.length > 0) {\n  generalMatches.items[0].insertText(\n    \"de la Universidad del Mar campus Puerto Escondido\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- 3) First \"Objetivo espec\u00edfico\" paragraph ---------------------------\n// \"...plantas existente del jard\u00edn bot\u00e1nico \"\u201d.\u201d  ->  \"...plantas\n// existente del jard\u00edn bot\u00e1nico de la Universidad del Mar.\" and the\n// \"_GoBack\" bookmark now sits right before the trailing period.\nconst specificMatches = body.search(\"\u201c\u201d\", { matchCase: false });\nspecificMatches.load(\"text\");\nawait context.sync();\nif (specificMatches.items.length > 0) {\n  specificMatches.items[0].insertText(\"de la Universidad del Mar\", \"Replace\");\n  await context.sync();\n}\n\n// Move the \"_GoBack\" bookmark to just before the final period of that\n// same paragraph (remove it from wherever it currently lives first).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (\n    text.indexOf(\"Objetivo espec\u00edfico\") === 0 &&\n    text.indexOf(\"existente del jard\u00edn bot\u00e1nico de la Universidad del Mar\") !==\n      -1\n  ) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const finalDot = targetParagraph.search(\".\", { matchCase: true });\n  finalDot.load(\"text\");\n  await context.sync();\n  if (finalDot.items.length > 0) {\n    const dotStart = finalDot.items[finalDot.items.length - 1].getRange(\n      \"Start\"\n    );\n    dotStart.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n\n// --- 4) Last paragraph: merge the two runs back into one ---------------\n// Re-typing the paragraph's own text over its whole range merges the\n// (now bookmark-free) runs back into a single run, matching the diff.\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst lastWhole = lastParagraph.getRange(\"Whole\");\nlastWhole.insertText(lastParagraph.text, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the empty curly-quote placeholders (\"\u201c\u201d\") throughout the\n# document with the real name of the botanical garden, and tidy up the\n# last paragraph (which previously held the stray \"_GoBack\" bookmark\n# used only to mark where the placeholder text used to be typed).\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph --------------------------------------------------\n# \"...jard\u00edn bot\u00e1nico \u201c\u201d.\"  ->  \"...jard\u00edn bot\u00e1nico de la Universidad del\n# Mar campus Puerto Escondido.\"\n$p1 = $d.Paragraphs(1).Range\n$null = $p1.Find.Execute(\"\u201c\u201d\", $false, $false, $false, $false, $false, $true, 1, $false, \"de la Universidad del Mar campus Puerto Escondido\", 1)\n\n# --- 2) \"Objetivo general\" paragraph -------------------------------------\n# \"...visitantes del jard\u00edn bot\u00e1nico \u201c\u201d de la Universidad del Mar.\"  ->\n# \"...visitantes del jard\u00edn bot\u00e1nico de la Universidad del Mar campus\n# Puerto Escondido.\"\n# Scope the find to the second run only (it starts at \". Dicha\n# aplicaci\u00f3n...\") so the untouched first run (\"Objetivo general: ...\n# Android\") is not merged into it.\n$p3 = $d.Paragraphs(3).Range\n$run2Start = $p3.Duplicate\n$null = $run2Start.Find.Execute(\". Dicha aplicaci\u00f3n\")\n$run2Range = $d.Range($run2Start.Start, $p3.End)\n$null = $run2Range.Find.Execute(\"\u201c\u201d de la Universidad del Mar\", $false, $false, $false, $false, $false, $true, 1, $false, \"de la Universidad del Mar campus Puerto Escondido\", 1)\n\n# --- 3) First \"Objetivo espec\u00edfico\" paragraph ----------------------------\n# \"...plantas existente del jard\u00edn bot\u00e1nico \u201c\u201d.\"  ->  \"...plantas\n# existente del jard\u00edn bot\u00e1nico de la Universidad del Mar.\" and the\n# \"_GoBack\" bookmark now sits right before the trailing period.\n$p4 = $d.Paragraphs(4).Range\n$null = $p4.Find.Execute(\"\u201c\u201d\", $false, $false, $false, $false, $false, $true, 1, $false, \"de la Universidad del Mar\", 1)\n\n# Move the \"_GoBack\" bookmark to just before the final period of that\n# same paragraph (remove it from wherever it currently lives first).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$p4 = $d.Paragraphs(4).Range\n$dotRange = $p4.Duplicate\n$null = $dotRange.Find.Execute(\".\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$bookmarkPoint = $d.Range($dotRange.Start, $dotRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n\n# --- 4) Last paragraph: merge the two runs back into one -----------------\n# Re-typing the paragraph's own text (excluding the trailing paragraph\n# mark) over its whole range merges the (now bookmark-free) runs back\n# into a single run, matching the diff.\n$paraCount = $d.Paragraphs.Count\n$lastParaFull = $d.Paragraphs($paraCount).Range\n$lastRange = $d.Range($lastParaFull.Start, $lastParaFull.End - 1)\n$lastText = $lastRange.Text\n$null = $lastRange.Find.Execute($lastText, $false, $false, $false, $false, $false, $true, 1, $false, $lastText, 1)\n"}
